$d = $word.ActiveDocument

# 1) Update the activation date.
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# Helper: find a paragraph whose trimmed text equals $text and return its 1-based index.
function Find-ParagraphIndex($doc, $text) {
    $idx = 0
    $foundIdx = -1
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            $foundIdx = $idx
        }
    }
    return $foundIdx
}

# Helper: insert a new paragraph right after paragraph #$afterIdx, set its text
# and make the run (but not the paragraph mark) italic.
function Insert-ItalicParagraphAfter($doc, $afterIdx, $newText) {
    $p = $doc.Paragraphs.Item($afterIdx)
    $p.Range.InsertParagraphAfter() | Out-Null
    $newP = $doc.Paragraphs.Item($afterIdx + 1)
    $newP.Range.Text = $newText
    $r = $newP.Range
    $r.End = $r.End - 1
    $r.Font.Italic = 1
}

# 2) Objetivos: add the English translation paragraph right after the Portuguese one.
$objIdx = Find-ParagraphIndex $d "Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais."
Insert-ItalicParagraphAfter $d $objIdx "To present the experimental techniques of materialographic preparation and characterization of materials."

# 3) Programa resumido: add the English translation paragraph right after the Portuguese one.
$resumidoIdx = Find-ParagraphIndex $d "Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica."
Insert-ItalicParagraphAfter $d $resumidoIdx "X-ray diffraction. Materialography. Optical microscopy. Electron microscopy. Thermal analysis."

# 4) Programa: merge the two runs (the manual line break in between becomes a single run/sentence).
$brChar = [string][char]11
$searchStr = "difração de raios X. " + $brChar + "Preparação materialográfica"
$d.Content.Find.Execute($searchStr, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "difração de raios X. Preparação materialográfica", 2) | Out-Null

# 5) Programa: add the English translation paragraph right after the (now merged) Portuguese one.
$programaText = "A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica."
$programaIdx = Find-ParagraphIndex $d $programaText
$englishProgramaText = "The microstructure of materials. Crystal lattices and systems, space groups and symmetry, most common types of crystal structures. Stereographic projection. Direction of the diffracted beam and Bragg's law. Intensity of the diffracted beam. Methods of X-ray diffraction.Materialographic sample preparation: cutting, embedding, sanding and polishing. Chemical etching techniques to reveal phases. Fundamentals of quantitative materialography. Optical microscopy. Electron microscopy techniques: scanning and transmission. Chemical analysis of microregions: energy dispersive spectroscopy. Thermal analysis techniques: differential thermal analysis, differential scanning calorimetry and thermogravimetric analysis."
Insert-ItalicParagraphAfter $d $programaIdx $englishProgramaText
